$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 94011938
$ws.Range("B2").Value = 98520
$ws.Range("E2").Value = 222498
$ws.Range("F2").Value = "Blåsippa"
$ws.Range("G2").Value = "Hepatica nobilis"
$ws.Range("H2").Value = "Schreb."
$ws.Range("P2").Value = "Ygne-Hemse, Gtl"
$ws.Range("Q2").Value = 698779.9316287825
$ws.Range("R2").Value = 6356965.133827931
$ws.Range("S2").Value = 10
$ws.Range("Y2").Value = "'2021-05-24"
$ws.Range("AA2").Value = "'2021-05-24"
$ws.Range("AW2").Value = "Tony Svensson"
$ws.Range("AX2").Value = "Tony Svensson"
$ws.Range("AY2").Value = "Ecogain"
$ws.Range("AC2").ClearContents() | Out-Null
$ws.Range("AH2").ClearContents() | Out-Null

# Row 3
$ws.Range("A3").Value = 94011937
$ws.Range("B3").Value = 98520
$ws.Range("E3").Value = 222498
$ws.Range("F3").Value = "Blåsippa"
$ws.Range("G3").Value = "Hepatica nobilis"
$ws.Range("H3").Value = "Schreb."
$ws.Range("P3").Value = "Ygne-Hemse, Gtl"
$ws.Range("Q3").Value = 698786.3346319427
$ws.Range("R3").Value = 6356944.887886292
$ws.Range("S3").Value = 10
$ws.Range("Y3").Value = "'2021-05-24"
$ws.Range("AA3").Value = "'2021-05-24"
$ws.Range("AW3").Value = "Tony Svensson"
$ws.Range("AX3").Value = "Tony Svensson"
$ws.Range("AY3").Value = "Ecogain"
$ws.Range("AC3").ClearContents() | Out-Null
$ws.Range("AH3").ClearContents() | Out-Null

# Row 4
$ws.Range("A4").Value = 94011889
$ws.Range("B4").Value = 56887
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 102995
$ws.Range("F4").Value = "Buskskvätta"
$ws.Range("G4").Value = "Saxicola rubetra"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("M4").Value = "spel/sång"
$ws.Range("P4").Value = "Ygne-Hemse, Gtl"
$ws.Range("Q4").Value = 698742.6195520113
$ws.Range("R4").Value = 6356963.322957435
$ws.Range("S4").Value = 10
$ws.Range("Y4").Value = "'2021-05-24"
$ws.Range("AA4").Value = "'2021-05-24"
$ws.Range("AW4").Value = "Tony Svensson"
$ws.Range("AX4").Value = "Tony Svensson"
$ws.Range("AY4").Value = "Ecogain"
$ws.Range("AC4").ClearContents() | Out-Null
$ws.Range("AH4").ClearContents() | Out-Null

# Row 5
$ws.Range("A5").Value = 2999289
$ws.Range("B5").Value = 96319
$ws.Range("E5").Value = 219799
$ws.Range("F5").Value = "Kärrknipprot"
$ws.Range("G5").Value = "Epipactis palustris"
$ws.Range("H5").Value = "(L.) Crantz"

# Row 6
$ws.Range("A6").Value = 5168355
$ws.Range("B6").Value = 97335
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 222662
$ws.Range("F6").Value = "Axag"
$ws.Range("G6").Value = "Schoenus ferrugineus"

# Row 7
$ws.Range("A7").Value = 2289844
$ws.Range("B7").Value = 96336
$ws.Range("E7").Value = 219811
$ws.Range("F7").Value = "Brudsporre"
$ws.Range("G7").Value = "Gymnadenia conopsea"
$ws.Range("H7").Value = "(L.) R. Br."
$ws.Range("P7").Value = "06J1A05, Gtl"
$ws.Range("Q7").Value = 698621.4564279296
$ws.Range("R7").Value = 6356902.266207782
$ws.Range("S7").Value = 50
$ws.Range("Y7").Value = "'2007-07-09"
$ws.Range("AA7").Value = "'2007-09-05"
$ws.Range("AC7").Value = "Lokalens storlek varierar starkt från 5 m i radie till flera hundra meter. Valde 50 m i denna rapportering.Mittkoordinater för varje rikkärrsobjekt, ej för varje observation."
$ws.Range("AH7").Value = "Rikkärr"
$ws.Range("AW7").Value = "Britta Johansson"
$ws.Range("AX7").Value = "Magnus Martinsson"
$ws.Range("AY7").Value = "Åtgärdsprogram för hotade arter"

# Row 8
$ws.Range("A8").Value = 4095111
$ws.Range("B8").Value = 95990
$ws.Range("E8").Value = 221930
$ws.Range("F8").Value = "Kärrlilja"
$ws.Range("G8").Value = "Tofieldia calyculata"
$ws.Range("H8").Value = "(L.) Wahlenb."
$ws.Range("P8").Value = "06J1A05, Gtl"
$ws.Range("Q8").Value = 698621.4564279296
$ws.Range("R8").Value = 6356902.266207782
$ws.Range("S8").Value = 50
$ws.Range("Y8").Value = "'2007-07-09"
$ws.Range("AA8").Value = "'2007-09-05"
$ws.Range("AC8").Value = "Lokalens storlek varierar starkt från 5 m i radie till flera hundra meter. Valde 50 m i denna rapportering.Mittkoordinater för varje rikkärrsobjekt, ej för varje observation."
$ws.Range("AH8").Value = "Rikkärr"
$ws.Range("AW8").Value = "Britta Johansson"
$ws.Range("AX8").Value = "Magnus Martinsson"
$ws.Range("AY8").Value = "Åtgärdsprogram för hotade arter"

# Row 9
$ws.Range("A9").Value = 3704166
$ws.Range("B9").Value = 103164
$ws.Range("E9").Value = 221137
$ws.Range("F9").Value = "Majviva"
$ws.Range("G9").Value = "Primula farinosa"
$ws.Range("H9").Value = "L."
$ws.Range("P9").Value = "06J1A05, Gtl"
$ws.Range("Q9").Value = 698621.4564279296
$ws.Range("R9").Value = 6356902.266207782
$ws.Range("S9").Value = 50
$ws.Range("Y9").Value = "'2007-07-09"
$ws.Range("AA9").Value = "'2007-09-05"
$ws.Range("AC9").Value = "Lokalens storlek varierar starkt från 5 m i radie till flera hundra meter. Valde 50 m i denna rapportering.Mittkoordinater för varje rikkärrsobjekt, ej för varje observation."
$ws.Range("AH9").Value = "Rikkärr"
$ws.Range("AW9").Value = "Britta Johansson"
$ws.Range("AX9").Value = "Magnus Martinsson"
$ws.Range("AY9").Value = "Åtgärdsprogram för hotade arter"
$ws.Range("M9").ClearContents() | Out-Null
